$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string (header) text updates ---
$volCell = $ws.Range("A8")
$volChars = $volCell.Characters(21, 2)
$volChars.Text = "51"

$weekCell = $ws.Range("C9")
$weekChars1 = $weekCell.Characters(27, 10)
$weekChars1.Text = "12/18/2023"
$weekChars2 = $weekCell.Characters(48, 10)
$weekChars2.Text = "12/24/2023"

# --- Crime-grid numeric updates ---
$numericUpdates = @{
  "C16" = 1
  "E16" = -50
  "F16" = 6
  "G16" = 5
  "H16" = 20
  "J16" = 91
  "K16" = 23.076923076923
  "L16" = 19.148936170212
  "N16" = -81.758957654723
  "C17" = 5
  "D17" = 2
  "E17" = 150
  "F17" = 8
  "G17" = 8
  "H17" = 0
  "I17" = 106
  "J17" = 114
  "K17" = -7.017543859649
  "L17" = -15.2
  "M17" = 100
  "N17" = -34.161490683229
  "D18" = 5
  "E18" = 20
  "G18" = 14
  "H18" = 28.571428571428
  "I18" = 132
  "J18" = 183
  "K18" = -27.868852459016
  "L18" = -10.204081632653
  "M18" = 38.947368421052
  "N18" = -89.099917423616
  "D19" = 8
  "E19" = 50
  "F19" = 56
  "G19" = 65
  "H19" = -13.846153846153
  "I19" = 677
  "J19" = 664
  "K19" = 1.957831325301
  "L19" = 25.836431226765
  "M19" = -8.016304347826
  "N19" = -69.924478009773
  "C20" = 1
  "E20" = -50
  "I20" = 62
  "J20" = 67
  "K20" = -7.462686567164
  "L20" = 34.782608695652
  "M20" = 106.666666666667
  "N20" = -90.127388535031
  "C21" = 25
  "D21" = 19
  "E21" = 31.578947368421
  "F21" = 101
  "H21" = 0
  "I21" = 1097
  "J21" = 1136
  "K21" = -3.433098591549
  "L21" = 14.509394572025
  "M21" = 12.167689161554
  "N21" = -77.515884402541
  "E22" = -100
  "G22" = 3
  "H22" = -66.666666666666
  "J22" = 32
  "K22" = 15.625
  "M22" = 2.777777777777
  "C24" = 21
  "D24" = 13
  "E24" = 61.538461538461
  "F24" = 78
  "G24" = 66
  "H24" = 18.181818181818
  "I24" = 1045
  "J24" = 1212
  "K24" = -13.778877887788
  "L24" = 5.876393110435
  "M24" = 68.820678513731
  "C25" = 3
  "D25" = 2
  "E25" = 50
  "F25" = 21
  "G25" = 21
  "H25" = 0
  "I25" = 246
  "J25" = 235
  "K25" = 4.680851063829
  "L25" = 35.911602209944
  "M25" = 9.821428571428
  "G27" = 2
  "H27" = -100
  "L30" = 57.142857142857
}
foreach ($addr in $numericUpdates.Keys) {
  $ws.Range($addr).Value = $numericUpdates[$addr]
}

# --- Cells that flip from numeric to the text markers "0" / "***.*" ---
# Copying from an existing cell that already carries the right shared-string + style (s=14)
# reproduces both the text value and the original formatting exactly.
$zeroSource = $ws.Range("C14")
$starSource = $ws.Range("L14")
$zeroSource.Copy($ws.Range("D14"))
$starSource.Copy($ws.Range("E14"))
$zeroSource.Copy($ws.Range("C22"))
$zeroSource.Copy($ws.Range("D27"))
$starSource.Copy($ws.Range("E27"))
$zeroSource.Copy($ws.Range("F27"))
